# Applies the commit:
#  1. Slide 5's table gets a new table style (GUID swap).
#  2. The presentation's theme colour scheme ("Integral" / "Red Violet")
#     is swapped for the stock "Office Theme" / "Office" colour scheme.
#
# Colours below are RGB hex values converted to the BGR-packed long that
# PowerPoint's COM RGBColor.RGB property expects (0x00BBGGRR).

$p = $ppt.ActivePresentation

# --- 1. Table style on slide 5, shape 2 (the graphicFrame/table) -----------
$slide5 = $p.Slides.Item(5)
$table  = $slide5.Shapes.Item(2).Table
$table.ApplyStyle("{736C3C3D-1838-4DDD-9DDD-D3135790314B}")

# --- 2. Swap the deck's colour theme from "Integral"/Red Violet to the ----
#        stock Office Theme colours.
$tcs = $p.Slides.Item(1).ThemeColorScheme

$tcs.Item(1).RGB  = 0x000000   # dk1      000000
$tcs.Item(2).RGB  = 0xFFFFFF   # lt1      FFFFFF
$tcs.Item(3).RGB  = 0x6A5444   # dk2      44546A
$tcs.Item(4).RGB  = 0xE6E6E7   # lt2      E7E6E6
$tcs.Item(5).RGB  = 0xD59B5B   # accent1  5B9BD5
$tcs.Item(6).RGB  = 0x317DED   # accent2  ED7D31
$tcs.Item(7).RGB  = 0xA5A5A5   # accent3  A5A5A5
$tcs.Item(8).RGB  = 0x00C0FF   # accent4  FFC000
$tcs.Item(9).RGB  = 0xC47244   # accent5  4472C4
$tcs.Item(10).RGB = 0x47AD70   # accent6  70AD47
$tcs.Item(11).RGB = 0xC16305   # hlink    0563C1
$tcs.Item(12).RGB = 0x724F95   # folHlink 954F72
